# A new weekly price record (row) is inserted before the existing row 88,
# pushing all subsequent rows (old 88..186) down by one (new 89..187).
# The newly inserted row 88 receives a fresh "Florida King" / "Primera"
# price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 88 (existing row 88 and everything
# below it shift down by one row).
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new data record.
$ws.Cells.Item(88, 1).Value = 7
$ws.Cells.Item(88, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(88, 3).Value = "Ñuble"
$ws.Cells.Item(88, 4).Value = 44539
$ws.Cells.Item(88, 5).Value = 16
$ws.Cells.Item(88, 6).Value = "Fruta"
$ws.Cells.Item(88, 7).Value = 100103
$ws.Cells.Item(88, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(88, 9).Value = 100103004
$ws.Cells.Item(88, 10).Value = "Durazno"
$ws.Cells.Item(88, 11).Value = "Florida King"
$ws.Cells.Item(88, 12).Value = "Primera"
$ws.Cells.Item(88, 13).Value = 300
$ws.Cells.Item(88, 14).Value = 16000
$ws.Cells.Item(88, 15).Value = 17000
$ws.Cells.Item(88, 16).Value = 16500
$ws.Cells.Item(88, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(88, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(88, 19).Value = 1031
$ws.Cells.Item(88, 20).Value = 16

# Make sure the date cell keeps the same date number format used by the
# rest of the "Fecha" column (style index 2 in the original workbook).
$ws.Cells.Item(88, 4).NumberFormat = $ws.Cells.Item(89, 4).NumberFormat
